# This script updates the "F" column (numeric view/attendee counts) values
# across the 展览 (Exhibition), 演出 (Show), and 全部类型 (All types) sheets,
# matching the regenerated gh-pages data snapshot at commit 456a3b4.
# The 本地生活 (Local life) sheet has no changes in this update.

$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 925
$ws.Range("F3").Value = 1027
$ws.Range("F4").Value = 813
$ws.Range("F5").Value = 888
$ws.Range("F6").Value = 469
$ws.Range("F7").Value = 720
$ws.Range("F8").Value = 170
$ws.Range("F9").Value = 1322
$ws.Range("F10").Value = 747
$ws.Range("F12").Value = 568
$ws.Range("F14").Value = 59
$ws.Range("F15").Value = 1223
$ws.Range("F16").Value = 149
$ws.Range("F18").Value = 434
$ws.Range("F22").Value = 165
$ws.Range("F24").Value = 39
$ws.Range("F25").Value = 1100

# ---- Sheet: 演出 ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 115
$ws.Range("F7").Value = 260
$ws.Range("F11").Value = 119

# ---- Sheet: 全部类型 ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 925
$ws.Range("F5").Value = 1027
$ws.Range("F6").Value = 813
$ws.Range("F7").Value = 888
$ws.Range("F8").Value = 469
$ws.Range("F9").Value = 469
$ws.Range("F10").Value = 720
$ws.Range("F11").Value = 170
$ws.Range("F12").Value = 1322
$ws.Range("F13").Value = 747
$ws.Range("F14").Value = 115
$ws.Range("F17").Value = 568
$ws.Range("F20").Value = 59
$ws.Range("F21").Value = 1223
$ws.Range("F23").Value = 149
$ws.Range("F25").Value = 434
$ws.Range("F28").Value = 260
$ws.Range("F33").Value = 119
$ws.Range("F34").Value = 119
$ws.Range("F35").Value = 165
$ws.Range("F37").Value = 39
$ws.Range("F38").Value = 1100

